# Insert a new row at position 53, shifting existing rows 53:103 down to 54:104.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("53:53").Insert(-4121)

# Populate the newly inserted row 53 with its data.
$ws.Range("A53").Value = 6
$ws.Range("B53").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C53").Value = 'Metropolitana'
$ws.Range("D53").Value = 44447
$ws.Range("E53").Value = 13
$ws.Range("F53").Value = 100112029
$ws.Range("G53").Value = 'Orégano'
$ws.Range("H53").Value = 'Sin especificar'
$ws.Range("I53").Value = 'Primera'
$ws.Range("J53").Value = 28
$ws.Range("K53").Value = 9000
$ws.Range("L53").Value = 10000
$ws.Range("M53").Value = 9464
$ws.Range("N53").Value = '$/docena de atados'
$ws.Range("O53").Value = 'Región Metropolitana'
$ws.Range("P53").Value = 3155
$ws.Range("Q53").Value = 3
$ws.Range("R53").Value = 'Hortaliza'
